$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoanDetails")

# Row 13: verifyDataValues
$ws.Range("A13").Value = "verifyDataValues"

# Row 14: header row again (CarLoanAmount / InterestRate / LoanTenure)
$ws.Range("A14").Value = "CarLoanAmount"
$ws.Range("B14").Value = "InterestRate"
$ws.Range("C14").Value = "LoanTenure"

# Row 15: qwerty/qwerty/qwerty
$ws.Range("A15").Value = "qwerty"
$ws.Range("B15").Value = "qwerty"
$ws.Range("C15").Value = "qwerty"

# Update selection to C20 as in diff
$ws.Range("C20").Select()
